$d = $word.ActiveDocument

$pairs = @(
    @("35×99=", "12×47="),
    @("41×27=", "96×93="),
    @("94×35=", "76×29="),
    @("68×54=", "26×83="),
    @("93×96=", "64×53="),
    @("87×79=", "52×78="),
    @("88×33=", "42×21="),
    @("35×77=", "21×38="),
    @("90×33=", "85×60="),
    @("26×89=", "27×91="),
    @("80×66=", "64×73="),
    @("13×65=", "63×95="),
    @("66×46=", "68×25="),
    @("41×82=", "52×37="),
    @("74×54=", "26×11="),
    @("86×74=", "26×52="),
    @("69×46=", "98×85="),
    @("65×42=", "34×40="),
    @("25×33=", "13×23="),
    @("30×78=", "82×94="),
    @("32×31=", "31×21="),
    @("81×40=", "92×79="),
    @("69×64=", "49×40="),
    @("29×71=", "81×80="),
    @("76×66=", "94×22=")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
